# Regenerate save_data to use K (strikeouts) instead of Strike# for the
# "K" column, and recompute dependent IP/IF figures for the affected row.
#
# The workbook stores per-appearance pitching data with columns:
#   A=idx  B=date  C=TB  D=PC  E=dS0  F=dSF  G=K  H=IP  I=I0  J=IF
#
# This edit rewrites the "K" (column G) values for the game log rows, and
# updates the IP/IF values (columns H/J) for row 5 to match the
# regenerated statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K (column G) values, keyed by spreadsheet row number.
$kValues = @{
    2  = 1
    3  = 2
    4  = 0
    5  = 1
    6  = 1
    7  = 0
    9  = 0
    10 = 1
    11 = 0
    12 = 0
    13 = 1
    14 = 0
    15 = 1
    16 = 1
    17 = 0
    18 = 1
    19 = 3
    20 = 1
    21 = 0
    22 = 0
    24 = 0
    26 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}

# Row 5 also regenerated IP (H) and IF (J) figures alongside K.
$ws.Range("H5").Value = 3
$ws.Range("J5").Value = 9
